$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Add ArticleTableViewCell"

$ws.Range("C10").Select()
